$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shapes on this slide (before edit), in z-order:
#   1 Title 1      (id=4)  - "Khulna University of Engineering & Technology / BECM 4000 / Project and Thesis"
#   2 Subtitle 2   (id=5)  - main heading ("MACHINE LEARNING ALGORITHMS ...")
#   3 Rectangle 2  (id=9)
#   4 Subtitle 2   (id=12) - "Presented By: ..."
#   5 Subtitle 2   (id=13) - "Supervised By: ..."
#   6 Picture 13   (id=14) - logo image

# Move the main heading box up now that the "Title 1" box above it is going away.
$s.Shapes.Item(2).Top = 164.0208

# Delete shapes from the highest index downward so earlier indices stay stable.
$s.Shapes.Item(6).Delete()   # Picture 13 (id=14)
$s.Shapes.Item(5).Delete()   # Subtitle 2 (id=13) "Supervised By"
$s.Shapes.Item(4).Delete()   # Subtitle 2 (id=12) "Presented By"
$s.Shapes.Item(1).Delete()   # Title 1 (id=4)
